# "Added confirmation-page and implemented logout-function (not tested yet)"
#
# - Mark the "Seite, die Übersicht über Warenkorb zeigt" (confirmation /
#   cart-overview page) TODO as done, highlighting it green.
# - Add two new open TODO items: sending a confirmation e-mail after an
#   order, and implementing the logout function.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark the confirmation-page TODO as done -------------------------------
$ws.Range("B2").Value2 = "done"
$ws.Range("B2").Font.ThemeColor = 2      # white font (Light1)
$ws.Range("B2").Interior.ThemeColor = 10 # green fill (Accent6)

# --- New TODO: send confirmation e-mail after an order ---------------------
$ws.Range("A9").Value2 = "E-Mail versenden nach Bestellung"
$ws.Range("B9").Value2 = "offen"
$ws.Range("B9").Style = "Schlecht"
$ws.Range("C9").Value2 = "Jonas"

# --- New TODO: implement logout -------------------------------------------
$ws.Range("A10").Value2 = "Logout"
$ws.Range("B10").Value2 = "offen"
$ws.Range("B10").Style = "Schlecht"
$ws.Range("C10").Value2 = "Jonas"

$ws.Range("C10").Select()
